$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H7").Value = 23000
$ws.Range("I7").Value = 15999
$ws.Range("K7").Value = 15999
$ws.Range("M7").Value = -15887
$ws.Range("H14").Value = 23000
$ws.Range("I14").Value = 15999
$ws.Range("K14").Value = 15999
$ws.Range("M14").Value = -15808
$ws.Range("H28").Value = 382.91666
$ws.Range("I28").Value = 382.91666
$ws.Range("K28").Value = 382.91666
$ws.Range("M28").Value = 102.08334
$ws.Range("H70").Value = 4411.25
$ws.Range("J70").Value = 4411.25
$ws.Range("L70").Value = 13233.75
$ws.Range("N70").Value = -13773.75
$ws.Range("H73").Value = 4411.25
$ws.Range("J73").Value = 4411.25
$ws.Range("L73").Value = 13233.75
$ws.Range("N73").Value = -15105.75
$ws.Range("H80").Value = 1099.625
$ws.Range("I80").Value = 1166.5
$ws.Range("J80").Value = 899
$ws.Range("K80").Value = 3499.5
$ws.Range("L80").Value = 2697
$ws.Range("M80").Value = -2501.5
$ws.Range("N80").Value = -4693
$ws.Range("H83").Value = 1099.625
$ws.Range("I83").Value = 1166.5
$ws.Range("J83").Value = 899
$ws.Range("K83").Value = 10498.5
$ws.Range("L83").Value = 8091
$ws.Range("M83").Value = -5506.5
$ws.Range("N83").Value = -18075
$ws.Range("H88").Value = 654.2
$ws.Range("I88").Value = 223
$ws.Range("J88").Value = 1301
$ws.Range("K88").Value = 223
$ws.Range("L88").Value = 1301
$ws.Range("M88").Value = 183
$ws.Range("N88").Value = -2113
$ws.Range("H91").Value = 654.2
$ws.Range("I91").Value = 223
$ws.Range("J91").Value = 1301
$ws.Range("K91").Value = 223
$ws.Range("L91").Value = 1301
$ws.Range("M91").Value = 1181
$ws.Range("N91").Value = -4109

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 0
$ws.Range("I45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("H61").Value = 1828
$ws.Range("I61").Value = 1828
$ws.Range("K61").Value = 1828
$ws.Range("M61").Value = -1616
$ws.Range("H63").Value = 2116.9
$ws.Range("I63").Value = 1294.8334
$ws.Range("J63").Value = 3350
$ws.Range("K63").Value = 1294.8334
$ws.Range("L63").Value = 3350
$ws.Range("M63").Value = -608.8334
$ws.Range("N63").Value = -4722
$ws.Range("H66").Value = 2116.9
$ws.Range("I66").Value = 1294.8334
$ws.Range("J66").Value = 3350
$ws.Range("K66").Value = 6474.166999999999
$ws.Range("L66").Value = 16750
$ws.Range("M66").Value = -3042.166999999999
$ws.Range("N66").Value = -23614
$ws.Range("H74").Value = 1624.5
$ws.Range("I74").Value = 1750
$ws.Range("J74").Value = 1499
$ws.Range("K74").Value = 1750
$ws.Range("L74").Value = 1499
$ws.Range("M74").Value = -876
$ws.Range("N74").Value = -3247
$ws.Range("H77").Value = 1624.5
$ws.Range("I77").Value = 1750
$ws.Range("J77").Value = 1499
$ws.Range("K77").Value = 8750
$ws.Range("L77").Value = 7495
$ws.Range("M77").Value = -4382
$ws.Range("N77").Value = -16231
$ws.Range("H88").Value = 2181.9524
$ws.Range("I88").Value = 671.5
$ws.Range("J88").Value = 2786.1333
$ws.Range("K88").Value = 671.5
$ws.Range("L88").Value = 2786.1333
$ws.Range("M88").Value = -265.5
$ws.Range("N88").Value = -3598.1333
$ws.Range("H91").Value = 2181.9524
$ws.Range("I91").Value = 671.5
$ws.Range("J91").Value = 2786.1333
$ws.Range("K91").Value = 671.5
$ws.Range("L91").Value = 2786.1333
$ws.Range("M91").Value = 732.5
$ws.Range("N91").Value = -5594.1333
$ws.Range("H110").Value = 4608.727
$ws.Range("I110").Value = 697
$ws.Range("K110").Value = 697
$ws.Range("M110").Value = 1348
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("H122").Value = 2622.5715
$ws.Range("I122").Value = 2618.6
$ws.Range("J122").Value = 2632.5
$ws.Range("K122").Value = 7855.799999999999
$ws.Range("L122").Value = 7897.5
$ws.Range("M122").Value = -5405.799999999999
$ws.Range("N122").Value = -12797.5
$ws.Range("H136").Value = 1828
$ws.Range("I136").Value = 1828
$ws.Range("K136").Value = 5484
$ws.Range("M136").Value = -2934
$ws.Range("M45").ClearContents()
$ws.Range("N119").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2499
$ws.Range("J20").Value = 2895.4
$ws.Range("L20").Value = 2895.4
$ws.Range("N20").Value = -3389.4
$ws.Range("H86").Value = 3239.5
$ws.Range("I86").Value = 2135.4285
$ws.Range("J86").Value = 4785.2
$ws.Range("K86").Value = 2135.4285
$ws.Range("L86").Value = 4785.2
$ws.Range("M86").Value = -1012.4285
$ws.Range("N86").Value = -7031.2
$ws.Range("H89").Value = 3239.5
$ws.Range("I89").Value = 2135.4285
$ws.Range("J89").Value = 4785.2
$ws.Range("K89").Value = 10677.1425
$ws.Range("L89").Value = 23926
$ws.Range("M89").Value = -5061.1425
$ws.Range("N89").Value = -35158
$ws.Range("H107").Value = 1250.1
$ws.Range("I107").Value = 1232
$ws.Range("J107").Value = 1413
$ws.Range("K107").Value = 1232
$ws.Range("L107").Value = 1413
$ws.Range("M107").Value = 688
$ws.Range("N107").Value = -5253

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 3351.25
$ws.Range("I62").Value = 3135
$ws.Range("K62").Value = 3135
$ws.Range("M62").Value = -2511
$ws.Range("H65").Value = 3351.25
$ws.Range("I65").Value = 3135
$ws.Range("K65").Value = 15675
$ws.Range("M65").Value = -12555
$ws.Range("H99").Value = 6081
$ws.Range("I99").Value = 6081
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 6081
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -4583
$ws.Range("H126").Value = 6081
$ws.Range("I126").Value = 6081
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 18243
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -15773
$ws.Range("N99").ClearContents()
$ws.Range("N126").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 40.142857
$ws.Range("I33").Value = 13.5
$ws.Range("J33").Value = 200
$ws.Range("K33").Value = 81
$ws.Range("L33").Value = 1200
$ws.Range("M33").Value = 202
$ws.Range("N33").Value = -1766
$ws.Range("H44").Value = 345.83334
$ws.Range("I44").Value = 296
$ws.Range("J44").Value = 352.0625
$ws.Range("K44").Value = 888
$ws.Range("L44").Value = 1056.1875
$ws.Range("M44").Value = -490
$ws.Range("N44").Value = -1852.1875
$ws.Range("H69").Value = 2033.3334
$ws.Range("I69").Value = 1200
$ws.Range("J69").Value = 2450
$ws.Range("K69").Value = 3600
$ws.Range("L69").Value = 7350
$ws.Range("M69").Value = -2789
$ws.Range("N69").Value = -8972
$ws.Range("H72").Value = 2033.3334
$ws.Range("I72").Value = 1200
$ws.Range("J72").Value = 2450
$ws.Range("K72").Value = 10800
$ws.Range("L72").Value = 22050
$ws.Range("M72").Value = -6744
$ws.Range("N72").Value = -30162
$ws.Range("H80").Value = 4900
$ws.Range("J80").Value = 4900
$ws.Range("L80").Value = 14700
$ws.Range("N80").Value = -16572
$ws.Range("H83").Value = 4900
$ws.Range("J83").Value = 4900
$ws.Range("L83").Value = 44100
$ws.Range("N83").Value = -53460
$ws.Range("H117").Value = 833.3333
$ws.Range("J117").Value = 1900
$ws.Range("L117").Value = 5700
$ws.Range("N117").Value = -12584

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 6948311.5
$ws.Range("I122").Value = 7816413
$ws.Range("K122").Value = 23449239
$ws.Range("M122").Value = -23446789
$ws.Range("H132").Value = 2000
$ws.Range("I132").Value = 2000
$ws.Range("K132").Value = 6000
$ws.Range("M132").Value = -3470

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 7127.4
$ws.Range("I40").Value = 7034.25
$ws.Range("K40").Value = 7034.25
$ws.Range("M40").Value = -6898.25
$ws.Range("H43").Value = 10014
$ws.Range("J43").Value = 10014
$ws.Range("L43").Value = 10014
$ws.Range("N43").Value = -10400
$ws.Range("H55").Value = 500
$ws.Range("I55").Value = 500
$ws.Range("K55").Value = 500
$ws.Range("M55").Value = -327
$ws.Range("H68").Value = 499.75
$ws.Range("I68").Value = 499.66666
$ws.Range("K68").Value = 499.66666
$ws.Range("M68").Value = 249.33334
$ws.Range("H71").Value = 499.75
$ws.Range("I71").Value = 499.66666
$ws.Range("K71").Value = 2498.3333
$ws.Range("M71").Value = 1245.6667
$ws.Range("H82").Value = 1451.6154
$ws.Range("I82").Value = 828.6667
$ws.Range("J82").Value = 1985.5714
$ws.Range("K82").Value = 828.6667
$ws.Range("L82").Value = 1985.5714
$ws.Range("M82").Value = -467.6667
$ws.Range("N82").Value = -2707.5714
$ws.Range("H85").Value = 1451.6154
$ws.Range("I85").Value = 828.6667
$ws.Range("J85").Value = 1985.5714
$ws.Range("K85").Value = 828.6667
$ws.Range("L85").Value = 1985.5714
$ws.Range("M85").Value = 419.3333
$ws.Range("N85").Value = -4481.5714
$ws.Range("H119").Value = 50000
$ws.Range("J119").Value = 50000
$ws.Range("L119").Value = 50000
$ws.Range("N119").Value = -59676
$ws.Range("H122").Value = 3430.1428
$ws.Range("I122").Value = 3430.1428
$ws.Range("K122").Value = 10290.4284
$ws.Range("M122").Value = -7840.428400000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 2462.3333
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("H65").Value = 2462.3333
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("H136").Value = 1549.591
$ws.Range("I136").Value = 1549.591
$ws.Range("K136").Value = 4648.772999999999
$ws.Range("M136").Value = -2098.772999999999
$ws.Range("N62").ClearContents()
$ws.Range("N65").ClearContents()
